$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.439.26"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.946.61"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.66"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.374"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("E10").Value = "  -8.13%  "
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("D12").Value = "2.232.18"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").Value = "1.949.09"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "36.410.11"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("E20").Value = "  -5.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "227.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.136"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0608"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.82%  "
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.47%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.41%  "
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -13.13%  "
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.358.32"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").Value = "2.123.31"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.99%  "
